$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the grade_type values (column D) first so the new shared strings
# "Auto Pass" and "Percentage" are registered ahead of the renamed IDs,
# matching the order produced by the original authoring session.
$ws.Range("D3").Value = "Auto Pass"
$ws.Range("D2").Value = "Percentage"

# Rename the IP ids in column B
$ws.Range("B2").Value = "Ren_IP_0005"
$ws.Range("B3").Value = "Ren_IP_0006"

# Move the active selection to B3, as reflected in the saved sheet view
$ws.Range("B3").Select()
